$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4668373067085838
$ws.Range("C2").Value = 0.03987753500712188
$ws.Range("D2").Value = 0.1873939752031788
$ws.Range("E2").Value = 0.166210065083412
$ws.Range("F2").Value = 1.484535831506911
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1890045645123664
$ws.Range("K2").Value = 0.4298630101824585
$ws.Range("N2").Value = 1.667300998228899
$ws.Range("O2").Value = 3.623486681340324

$ws.Range("B3").Value = 0.4287475654600996
$ws.Range("C3").Value = 0.03478104256781478
$ws.Range("D3").Value = 0.1815751677574582
$ws.Range("E3").Value = 0.1620237846593824
$ws.Range("F3").Value = 1.485436559375891
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1849979845659249
$ws.Range("K3").Value = 0.3893307655360161
$ws.Range("N3").Value = 1.68493773425644
$ws.Range("O3").Value = 3.638125709781946

$ws.Range("B4").Value = 0.4054897257898915
$ws.Range("C4").Value = 0.03164125845404442
$ws.Range("D4").Value = 0.1780815867067247
$ws.Range("E4").Value = 0.1595368129236157
$ws.Range("F4").Value = 1.486746116947138
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1826434735646245
$ws.Range("K4").Value = 0.3645219074927866
$ws.Range("N4").Value = 1.696319255513855
$ws.Range("O4").Value = 3.648967928411309

$ws.Range("B5").Value = 0.3960450350762699
$ws.Range("C5").Value = 0.03035919039130874
$ws.Range("D5").Value = 0.1766779432694818
$ws.Range("E5").Value = 0.1585443799780499
$ws.Range("F5").Value = 1.487470109109651
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1817105536293866
$ws.Range("K5").Value = 0.3544322752698292
$ws.Range("N5").Value = 1.701096241540262
$ws.Range("O5").Value = 3.653852443779769

$ws.Range("B6").Value = 0.3944787640630807
$ws.Range("C6").Value = 0.03014614999634091
$ws.Range("D6").Value = 0.1764460815806927
$ws.Range("E6").Value = 0.1583808588062503
$ws.Range("F6").Value = 1.487601826499009
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1815572479465288
$ws.Range("K6").Value = 0.3527581341586483
$ws.Range("N6").Value = 1.701897846966538
$ws.Range("O6").Value = 3.6546916747252

$ws.Range("B7").Value = 0.4053622168021889
$ws.Range("C7").Value = 0.03162397839862763
$ws.Range("D7").Value = 0.1780625754816612
$ws.Range("E7").Value = 0.1595233434150174
$ws.Range("F7").Value = 1.4867551101513
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1826307843000876
$ws.Range("K7").Value = 0.3643857527664238
$ws.Range("N7").Value = 1.696383117128856
$ws.Range("O7").Value = 3.649031914882954

$ws.Range("B8").Value = 0.4536774065601321
$ws.Range("C8").Value = 0.03812248416735997
$ws.Range("D8").Value = 0.1853712760128872
$ws.Range("E8").Value = 0.1647493505914497
$ws.Range("F8").Value = 1.484689449904387
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1876012014551662
$ws.Range("K8").Value = 0.4158715922708325
$ws.Range("N8").Value = 1.673267472992481
$ws.Range("O8").Value = 3.628149500113437

$ws.Range("B9").Value = 0.549432037787227
$ws.Range("C9").Value = 0.05078036546186127
$ws.Range("D9").Value = 0.2003281708281008
$ws.Range("E9").Value = 0.1756578958779897
$ws.Range("F9").Value = 1.486638268594106
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.198185394669764
$ws.Range("K9").Value = 0.5174363035354759
$ws.Range("N9").Value = 1.632321110503866
$ws.Range("O9").Value = 3.601907847604906

$ws.Range("B10").Value = 0.6203807254739218
$ws.Range("C10").Value = 0.06002577945309895
$ws.Range("D10").Value = 0.2116939613104734
$ws.Range("E10").Value = 0.1840739315295679
$ws.Range("F10").Value = 1.491726132786368
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.2064727756213216
$ws.Range("K10").Value = 0.5924049034930761
$ws.Range("N10").Value = 1.604908412260047
$ws.Range("O10").Value = 3.591598606151763

$ws.Range("B11").Value = 0.6527837088962656
$ws.Range("C11").Value = 0.06421956658070371
$ws.Range("D11").Value = 0.2169456804096939
$ws.Range("E11").Value = 0.1879896311423863
$ws.Range("F11").Value = 1.494834492637878
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.2103541482696585
$ws.Range("K11").Value = 0.6265825886493701
$ws.Range("N11").Value = 1.593017225745978
$ws.Range("O11").Value = 3.588857641904866

$ws.Range("B12").Value = 0.6650718548515329
$ws.Range("C12").Value = 0.0658058658539602
$ws.Range("D12").Value = 0.2189459793983417
$ws.Range("E12").Value = 0.189484907157258
$ws.Range("F12").Value = 1.496125647310578
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2118399368558244
$ws.Range("K12").Value = 0.6395349956353584
$ws.Range("N12").Value = 1.588597615513486
$ws.Range("O12").Value = 3.588099970784015

$ws.Range("B13").Value = 0.6624245985358073
$ws.Range("C13").Value = 0.06546430898012545
$ws.Range("D13").Value = 0.2185146656922399
$ws.Range("E13").Value = 0.1891623182800686
$ws.Range("F13").Value = 1.495842501159558
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.2115192345898578
$ws.Range("K13").Value = 0.6367450251273681
$ws.Range("N13").Value = 1.589545750708303
$ws.Range("O13").Value = 3.588250682931516

$ws.Range("B14").Value = 0.6537943078372734
$ws.Range("C14").Value = 0.06435010879852143
$ws.Range("D14").Value = 0.2171100146396299
$ws.Range("E14").Value = 0.1881123984819268
$ws.Range("F14").Value = 1.494938430656603
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.2104760645166692
$ws.Range("K14").Value = 0.6276479927724949
$ws.Range("N14").Value = 1.592651951424489
$ws.Range("O14").Value = 3.588789690519292

$ws.Range("B15").Value = 0.6485103116744995
$ws.Range("C15").Value = 0.06366739288363021
$ws.Range("D15").Value = 0.2162511309136477
$ws.Range("E15").Value = 0.1874709166665198
$ws.Range("F15").Value = 1.494399516668324
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2098391753365263
$ws.Range("K15").Value = 0.622077090143847
$ws.Range("N15").Value = 1.594565442318435
$ws.Range("O15").Value = 3.589156349269615

$ws.Range("B16").Value = 0.6182656420854471
$ws.Range("C16").Value = 0.05975145804789861
$ws.Range("D16").Value = 0.2113523767057188
$ws.Range("E16").Value = 0.183819781476366
$ws.Range("F16").Value = 1.491538962727219
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2062213585675465
$ws.Range("K16").Value = 0.5901727514635127
$ws.Range("N16").Value = 1.605697185577163
$ws.Range("O16").Value = 3.591816952347699

$ws.Range("B17").Value = 0.5997439082712788
$ws.Range("C17").Value = 0.05734603491167434
$ws.Range("D17").Value = 0.2083679085376104
$ws.Range("E17").Value = 0.1816022246264382
$ws.Range("F17").Value = 1.489987383283292
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2040304619399507
$ws.Range("K17").Value = 0.5706190266997169
$ws.Range("N17").Value = 1.612674488373195
$ws.Range("O17").Value = 3.593948301275077

$ws.Range("B18").Value = 0.58910278588192
$ws.Range("C18").Value = 0.05596137410451263
$ws.Range("D18").Value = 0.2066589865261648
$ws.Range("E18").Value = 0.1803349557725227
$ws.Range("F18").Value = 1.489169682377991
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.2027808031475899
$ws.Range("K18").Value = 0.5593792509650939
$ws.Range("N18").Value = 1.616742141277424
$ws.Range("O18").Value = 3.595357620847835

$ws.Range("B19").Value = 0.585501978552486
$ws.Range("C19").Value = 0.05549236070318386
$ws.Range("D19").Value = 0.2060816948337276
$ws.Range("E19").Value = 0.1799072922367415
$ws.Range("F19").Value = 1.488905659406811
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.2023594919804168
$ws.Range("K19").Value = 0.555574878058394
$ws.Range("N19").Value = 1.61812873589205
$ws.Range("O19").Value = 3.595866294196583

$ws.Range("B20").Value = 0.6017143314641942
$ws.Range("C20").Value = 0.05760221326718806
$ws.Range("D20").Value = 0.2086848177753069
$ws.Range("E20").Value = 0.181837437902459
$ws.Range("F20").Value = 1.490144818154789
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.2042626014225135
$ws.Range("K20").Value = 0.5726998329275261
$ws.Range("N20").Value = 1.61192610229706
$ws.Range("O20").Value = 3.593702431526907

$ws.Range("B21").Value = 0.6563287538604925
$ws.Range("C21").Value = 0.06467742573737212
$ws.Range("D21").Value = 0.2175222811401625
$ws.Range("E21").Value = 0.1884204469554192
$ws.Range("F21").Value = 1.495200882373851
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2107820348802392
$ws.Range("K21").Value = 0.6303197434259573
$ws.Range("N21").Value = 1.591737323045715
$ws.Range("O21").Value = 3.588623764064607

$ws.Range("B22").Value = 0.6921261495433555
$ws.Range("C22").Value = 0.06929097923676864
$ws.Range("D22").Value = 0.2233655592523149
$ws.Range("E22").Value = 0.1927955739554008
$ws.Range("F22").Value = 1.499170261686686
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2151360936170761
$ws.Range("K22").Value = 0.6680360431938084
$ws.Range("N22").Value = 1.579028497734905
$ws.Range("O22").Value = 3.586938224310785

$ws.Range("B23").Value = 0.6730111066829068
$ws.Range("C23").Value = 0.0668296239200572
$ws.Range("D23").Value = 0.220240754304541
$ws.Range("E23").Value = 0.1904538476578281
$ws.Range("F23").Value = 1.496990909749528
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.2128037272822638
$ws.Range("K23").Value = 0.6479010011785533
$ws.Range("N23").Value = 1.58576697192245
$ws.Range("O23").Value = 3.587688334067963

$ws.Range("B24").Value = 0.6008234804529877
$ws.Range("C24").Value = 0.05748640048216203
$ws.Range("D24").Value = 0.2085415216443636
$ws.Range("E24").Value = 0.1817310741976925
$ws.Range("F24").Value = 1.490073410332286
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2041576202719142
$ws.Range("K24").Value = 0.5717590944036033
$ws.Range("N24").Value = 1.612264272245483
$ws.Range("O24").Value = 3.593813016155508

$ws.Range("B25").Value = 0.5234216593057113
$ws.Range("C25").Value = 0.04736546256532392
$ws.Range("D25").Value = 0.1962154622044636
$ws.Range("E25").Value = 0.172636256314739
$ws.Range("F25").Value = 1.485468914126983
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1952323776748983
$ws.Range("K25").Value = 0.4898978073259173
$ws.Range("N25").Value = 1.64292909718014
$ws.Range("O25").Value = 3.60743186500261

